$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header for column E ("URL" -> "Resource Description")
$ws.Range("E1").Value = "Resource Description"

# Add new "Description" values alongside existing "Subtitle" values
$ws.Range("E2").Value = "Description 1"
$ws.Range("E3").Value = "Description 2"
$ws.Range("E5").Value = "Description 4"

# Add new row 6 with B6 = 5
$ws.Range("B6").Value = 5
$ws.Rows.Item(6).RowHeight = 14.65

# Update the active selection to match the target state
$ws.Range("E6").Select()
